# VTQaZ (Vehicle Technologies Qualifying as ZEVs) update
#
# Commit: "Adds FF55 scenario with road CO2 standard effects integrated;
# uses FoPITY-2. Assumes ... H2 vehicles do not qualify as ZEVs (to align
# with ICCT HDV study)"
#
# Content change: row 8 ("hydrogen vehicle") on the VTQaZ sheet flips from
# 1 (qualifies as ZEV) to 0 (no longer qualifies) for every year column
# (B:AF, i.e. 2020-2050).

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("VTQaZ")

# Hydrogen vehicle is row 8 (column A label "hydrogen vehicle"); set every
# year value (columns B through AF) to 0.
$hydrogenRow = $wsData.Range("B8:AF8")
$hydrogenRow.Value = 0

# Reflect the author's last on-screen state: the VTQaZ sheet tab active
# with cell D6 selected (matches the saved sheetView/selection).
$wsData.Activate()
$wsData.Range("D6").Select()
